$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value2 = 4174.857
$ws.Range("I20").Value2 = 704
$ws.Range("K20").Value2 = 704
$ws.Range("M20").Value2 = -474

$ws.Range("H32").Value2 = 2000.5
$ws.Range("I32").Value2 = 2000.5
$ws.Range("K32").Value2 = 2000.5
$ws.Range("M32").Value2 = -1674.5

$ws.Range("H34").Value2 = 13944
$ws.Range("I34").Value2 = 13944
$ws.Range("J34").Value2 = 0
$ws.Range("K34").Value2 = 13944
$ws.Range("L34").Value2 = 0
$ws.Range("N34").ClearContents()
$ws.Range("M34").Value2 = -13741

$ws.Range("H35").Value2 = 4174.857
$ws.Range("I35").Value2 = 704
$ws.Range("K35").Value2 = 704
$ws.Range("M35").Value2 = -325

$ws.Range("H36").Value2 = 13944
$ws.Range("I36").Value2 = 13944
$ws.Range("J36").Value2 = 0
$ws.Range("K36").Value2 = 13944
$ws.Range("L36").Value2 = 0
$ws.Range("N36").ClearContents()
$ws.Range("M36").Value2 = -13229

$ws.Range("H40").Value2 = 2833.1667
$ws.Range("I40").Value2 = 3999.5
$ws.Range("J40").Value2 = 2250
$ws.Range("K40").Value2 = 3999.5
$ws.Range("L40").Value2 = 2250
$ws.Range("M40").Value2 = -3824.5
$ws.Range("N40").Value2 = -2600

$ws.Range("H51").Value2 = 10285.714
$ws.Range("I51").Value2 = 4000
$ws.Range("J51").Value2 = 15000
$ws.Range("K51").Value2 = 4000
$ws.Range("L51").Value2 = 15000
$ws.Range("M51").Value2 = -3516
$ws.Range("N51").Value2 = -15968

$ws.Range("H54").Value2 = 8420
$ws.Range("J54").Value2 = 8420
$ws.Range("L54").Value2 = 8420
$ws.Range("N54").Value2 = -9392

$ws.Range("H74").Value2 = 3833.3333
$ws.Range("J74").Value2 = 4250
$ws.Range("L74").Value2 = 4250
$ws.Range("N74").Value2 = -6122

$ws.Range("H77").Value2 = 3833.3333
$ws.Range("J77").Value2 = 4250
$ws.Range("L77").Value2 = 21250
$ws.Range("N77").Value2 = -30610

$ws.Range("H138").Value2 = 3838.2856
$ws.Range("J138").Value2 = 3630.2
$ws.Range("L138").Value2 = 10890.6
$ws.Range("N138").Value2 = -21170.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 964.6667
$ws.Range("I2").Value2 = 964.6667
$ws.Range("K2").Value2 = 964.6667
$ws.Range("M2").Value2 = -851.6667

$ws.Range("H32").Value2 = 4379.2
$ws.Range("I32").Value2 = 4379.2
$ws.Range("K32").Value2 = 4379.2
$ws.Range("M32").Value2 = -4092.2

$ws.Range("H34").Value2 = 29999
$ws.Range("I34").Value2 = 0
$ws.Range("J34").Value2 = 29999
$ws.Range("K34").Value2 = 0
$ws.Range("L34").Value2 = 29999
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value2 = -30541

$ws.Range("H61").Value2 = 5250
$ws.Range("I61").Value2 = 5250
$ws.Range("K61").Value2 = 5250
$ws.Range("M61").Value2 = -5038

$ws.Range("H74").Value2 = 0
$ws.Range("I74").Value2 = 0
$ws.Range("J74").Value2 = 0
$ws.Range("K74").Value2 = 0
$ws.Range("L74").Value2 = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value2 = 0
$ws.Range("I77").Value2 = 0
$ws.Range("J77").Value2 = 0
$ws.Range("K77").Value2 = 0
$ws.Range("L77").Value2 = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

$ws.Range("H97").Value2 = 10155.5
$ws.Range("I97").Value2 = 218.33333
$ws.Range("K97").Value2 = 218.33333
$ws.Range("M97").Value2 = 277.66667

$ws.Range("H116").Value2 = 964.6667
$ws.Range("I116").Value2 = 964.6667
$ws.Range("K116").Value2 = 964.6667
$ws.Range("M116").Value2 = 1329.3333

$ws.Range("H136").Value2 = 5250
$ws.Range("I136").Value2 = 5250
$ws.Range("K136").Value2 = 15750
$ws.Range("M136").Value2 = -13200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 964.6667
$ws.Range("I3").Value2 = 964.6667
$ws.Range("K3").Value2 = 964.6667
$ws.Range("M3").Value2 = -850.6667

$ws.Range("H22").Value2 = 361.8
$ws.Range("I22").Value2 = 411.57144
$ws.Range("J22").Value2 = 245.66667
$ws.Range("K22").Value2 = 411.57144
$ws.Range("L22").Value2 = 245.66667
$ws.Range("M22").Value2 = -238.57144
$ws.Range("N22").Value2 = -591.6666700000001

$ws.Range("H38").Value2 = 1
$ws.Range("I38").Value2 = 1
$ws.Range("K38").Value2 = 1
$ws.Range("M38").Value2 = 415

$ws.Range("H105").Value2 = 4135.2
$ws.Range("I105").Value2 = 4135.2
$ws.Range("K105").Value2 = 4135.2
$ws.Range("M105").Value2 = -2388.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 4448.8335
$ws.Range("I31").Value2 = 3538.8
$ws.Range("J31").Value2 = 8999
$ws.Range("K31").Value2 = 3538.8
$ws.Range("L31").Value2 = 8999
$ws.Range("M31").Value2 = -3243.8
$ws.Range("N31").Value2 = -9589

$ws.Range("H34").Value2 = 4448.8335
$ws.Range("I34").Value2 = 3538.8
$ws.Range("J34").Value2 = 8999
$ws.Range("K34").Value2 = 3538.8
$ws.Range("L34").Value2 = 8999
$ws.Range("M34").Value2 = -3336.8
$ws.Range("N34").Value2 = -9403

$ws.Range("H35").Value2 = 693.8889
$ws.Range("I35").Value2 = 693.8889
$ws.Range("J35").Value2 = 0
$ws.Range("K35").Value2 = 693.8889
$ws.Range("L35").Value2 = 0
$ws.Range("M35").Value2 = -399.8889
$ws.Range("N35").ClearContents()

$ws.Range("H55").Value2 = 8073
$ws.Range("I55").Value2 = 8073
$ws.Range("K55").Value2 = 8073
$ws.Range("M55").Value2 = -7758

$ws.Range("H59").Value2 = 43999
$ws.Range("I59").Value2 = 41998
$ws.Range("J59").Value2 = 44666
$ws.Range("K59").Value2 = 41998
$ws.Range("L59").Value2 = 44666
$ws.Range("M59").Value2 = -40853
$ws.Range("N59").Value2 = -46956

$ws.Range("H60").Value2 = 30333
$ws.Range("I60").Value2 = 0
$ws.Range("J60").Value2 = 30333
$ws.Range("K60").Value2 = 0
$ws.Range("L60").Value2 = 30333
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value2 = -31355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value2 = 82.666664
$ws.Range("I47").Value2 = 82.666664
$ws.Range("K47").Value2 = 247.999992
$ws.Range("M47").Value2 = 183.000008

$ws.Range("H68").Value2 = 0
$ws.Range("J68").Value2 = 0
$ws.Range("L68").Value2 = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value2 = 0
$ws.Range("J71").Value2 = 0
$ws.Range("L71").Value2 = 0
$ws.Range("N71").ClearContents()

$ws.Range("H131").Value2 = 922.0833
$ws.Range("J131").Value2 = 945.36365
$ws.Range("L131").Value2 = 2836.09095
$ws.Range("N131").Value2 = -12916.09095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value2 = 11660.125
$ws.Range("I46").Value2 = 4466.6665
$ws.Range("J46").Value2 = 15976.2
$ws.Range("K46").Value2 = 4466.6665
$ws.Range("L46").Value2 = 15976.2
$ws.Range("M46").Value2 = -4310.6665
$ws.Range("N46").Value2 = -16288.2

$ws.Range("H62").Value2 = 45000
$ws.Range("I62").Value2 = 45000
$ws.Range("K62").Value2 = 45000
$ws.Range("M62").Value2 = -44314

$ws.Range("H65").Value2 = 45000
$ws.Range("I65").Value2 = 45000
$ws.Range("K65").Value2 = 135000
$ws.Range("M65").Value2 = -131568

$ws.Range("H80").Value2 = 1000
$ws.Range("I80").Value2 = 1000
$ws.Range("J80").Value2 = 1000
$ws.Range("K80").Value2 = 1000
$ws.Range("L80").Value2 = 1000
$ws.Range("M80").Value2 = -2
$ws.Range("N80").Value2 = -2996

$ws.Range("H83").Value2 = 1000
$ws.Range("I83").Value2 = 1000
$ws.Range("J83").Value2 = 1000
$ws.Range("K83").Value2 = 5000
$ws.Range("L83").Value2 = 5000
$ws.Range("M83").Value2 = -8
$ws.Range("N83").Value2 = -14984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value2 = 15000
$ws.Range("I20").Value2 = 15000
$ws.Range("J20").Value2 = 15000
$ws.Range("K20").Value2 = 15000
$ws.Range("L20").Value2 = 15000
$ws.Range("M20").Value2 = -14774
$ws.Range("N20").Value2 = -15452

$ws.Range("H46").Value2 = 1565.4736
$ws.Range("I46").Value2 = 1249.2858
$ws.Range("J46").Value2 = 1749.9166
$ws.Range("K46").Value2 = 1249.2858
$ws.Range("L46").Value2 = 1749.9166
$ws.Range("M46").Value2 = -1061.2858
$ws.Range("N46").Value2 = -2125.9166

$ws.Range("H55").Value2 = 607.1739
$ws.Range("I55").Value2 = 139.375
$ws.Range("K55").Value2 = 139.375
$ws.Range("M55").Value2 = 33.625

$ws.Range("H68").Value2 = 1522.7142
$ws.Range("I68").Value2 = 922.25
$ws.Range("J68").Value2 = 2323.3333
$ws.Range("K68").Value2 = 922.25
$ws.Range("L68").Value2 = 2323.3333
$ws.Range("M68").Value2 = -173.25
$ws.Range("N68").Value2 = -3821.3333

$ws.Range("H71").Value2 = 1522.7142
$ws.Range("I71").Value2 = 922.25
$ws.Range("J71").Value2 = 2323.3333
$ws.Range("K71").Value2 = 4611.25
$ws.Range("L71").Value2 = 11616.6665
$ws.Range("M71").Value2 = -867.25
$ws.Range("N71").Value2 = -19104.6665

$ws.Range("H82").Value2 = 3299.4707
$ws.Range("I82").Value2 = 849
$ws.Range("J82").Value2 = 4320.5
$ws.Range("K82").Value2 = 849
$ws.Range("L82").Value2 = 4320.5
$ws.Range("M82").Value2 = -488
$ws.Range("N82").Value2 = -5042.5

$ws.Range("H85").Value2 = 3299.4707
$ws.Range("I85").Value2 = 849
$ws.Range("J85").Value2 = 4320.5
$ws.Range("K85").Value2 = 849
$ws.Range("L85").Value2 = 4320.5
$ws.Range("M85").Value2 = 399
$ws.Range("N85").Value2 = -6816.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value2 = 4648099.5
$ws.Range("I100").Value2 = 13940120
$ws.Range("K100").Value2 = 27880240
$ws.Range("M100").Value2 = -27879699

$ws.Range("H117").Value2 = 41998.5
$ws.Range("J117").Value2 = 41998.5
$ws.Range("L117").Value2 = 41998.5
$ws.Range("N117").Value2 = -51176.5

$ws.Range("H135").Value2 = 0
$ws.Range("J135").Value2 = 0
$ws.Range("L135").Value2 = 0
$ws.Range("N135").ClearContents()

